$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three cells whose text content changed.
$ws.Range("D8").Value = "[{Status=Success},`n{username=user2, alias=}]"
$ws.Range("D9").Value = "[{Status=Success}, {Username=user_1, alias=},{Username=user_2, alias=blabla},…]"
$ws.Range("E8").Value = 'Add a contact (adding both directions), and alias=""'

# Update the view state: scroll so A5 is the top-left cell, and select E9.
$ws.Range("E9").Select()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
